$wb = $excel.ActiveWorkbook

# --- "cat" worksheet: row 7 now references the formbuilder worksheet via
#     an "include" record instead of a hardcoded git pattern ---
$cat = $wb.Worksheets.Item("cat")
$cat.Range("A7").Value = "include"
$cat.Range("B7").Value = "formbuilder"
$cat.Range("C7").ClearContents()
$cat.Range("D7").ClearContents()
$cat.PageSetup.Orientation = 2
$cat.Range("A7:D7").Select()

# --- "cim" worksheet: same row-level change, on row 4 ---
$cim = $wb.Worksheets.Item("cim")
$cim.Range("A4").Value = "include"
$cim.Range("B4").Value = "formbuilder"
$cim.Range("C4").ClearContents()
$cim.Range("D4").ClearContents()
$cim.PageSetup.Orientation = 2

# --- "pdfgen" worksheet: gains an extra blank/separator row and becomes
#     the active sheet/tab ---
$pdfgen = $wb.Worksheets.Item("pdfgen")
$cim.Range("A13:D13").Copy()
$pdfgen.Range("A8:D8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$pdfgen.Activate()
$pdfgen.Range("B12").Select()

# --- "formbuilder" worksheet: selection moves to the blank row below its
#     single data row ---
$formbuilder = $wb.Worksheets.Item("formbuilder")
$formbuilder.Range("A4:D4").Select()

$pdfgen.Activate()
